# Updates Efnb3-Ephb6 LR-pair sheet with recomputed TPM-derived NATMI
# statistics: the 'Neutrophils' sending-cluster block is replaced by a new
# 'ECs' sending-cluster block, every other block's expression-derived
# metrics (columns E:T) are recalculated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = 'ECs'
$ws.Range("G2").Value = 0.2163973333333333
$ws.Range("H2").Value = 0.649192
$ws.Range("I2").Value = 0.2690509661419307
$ws.Range("J2").Value = 0.2690509661419307
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1474273333333333
$ws.Range("N2").Value = 0.442282
$ws.Range("O2").Value = 0.1588601259223368
$ws.Range("P2").Value = 0.1588601259223368
$ws.Range("Q2").Value = 0.03190288179377777
$ws.Range("R2").Value = 0.287125936144
$ws.Range("S2").Value = 0.04274147036083347
$ws.Range("T2").Value = 0.04274147036083347
# Row 3
$ws.Range("A3").Value = 'ECs'
$ws.Range("G3").Value = 0.2163973333333333
$ws.Range("H3").Value = 0.649192
$ws.Range("I3").Value = 0.2690509661419307
$ws.Range("J3").Value = 0.2690509661419307
$ws.Range("O3").Value = 0.4626735347223893
$ws.Range("P3").Value = 0.4626735347223893
$ws.Range("Q3").Value = 0.09291582139733333
$ws.Range("R3").Value = 0.836242392576
$ws.Range("S3").Value = 0.124482761525361
$ws.Range("T3").Value = 0.124482761525361
# Row 4
$ws.Range("A4").Value = 'ECs'
$ws.Range("G4").Value = 0.2163973333333333
$ws.Range("H4").Value = 0.649192
$ws.Range("I4").Value = 0.2690509661419307
$ws.Range("J4").Value = 0.2690509661419307
$ws.Range("M4").Value = 0.01780266666666666
$ws.Range("N4").Value = 0.053408
$ws.Range("O4").Value = 0.01918323966442261
$ws.Range("P4").Value = 0.01918323966442261
$ws.Range("Q4").Value = 0.003852449592888889
$ws.Range("R4").Value = 0.034672046336
$ws.Range("S4").Value = 0.005161269165445109
$ws.Range("T4").Value = 0.00516126916544511
# Row 5
$ws.Range("A5").Value = 'ECs'
$ws.Range("G5").Value = 0.2163973333333333
$ws.Range("H5").Value = 0.649192
$ws.Range("I5").Value = 0.2690509661419307
$ws.Range("J5").Value = 0.2690509661419307
$ws.Range("M5").Value = 0.3334263333333333
$ws.Range("N5").Value = 1.000279
$ws.Range("O5").Value = 0.3592830996908513
$ws.Range("P5").Value = 0.3592830996908513
$ws.Range("Q5").Value = 0.07215256939644445
$ws.Range("R5").Value = 0.6493731245679999
$ws.Range("S5").Value = 0.09666546509029113
$ws.Range("T5").Value = 0.09666546509029113
# Row 6
$ws.Range("A6").Value = 'Inflammatory-Mac'
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.06315533333333333
$ws.Range("H6").Value = 0.189466
$ws.Range("I6").Value = 0.07852224049441002
$ws.Range("J6").Value = 0.07852224049441003
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1474273333333333
$ws.Range("N6").Value = 0.442282
$ws.Range("O6").Value = 0.1588601259223368
$ws.Range("P6").Value = 0.1588601259223368
$ws.Range("Q6").Value = 0.00931082237911111
$ws.Range("R6").Value = 0.083797401412
$ws.Range("S6").Value = 0.01247405301264599
$ws.Range("T6").Value = 0.01247405301264599
# Row 7
$ws.Range("A7").Value = 'Inflammatory-Mac'
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.06315533333333333
$ws.Range("H7").Value = 0.189466
$ws.Range("I7").Value = 0.07852224049441002
$ws.Range("J7").Value = 0.07852224049441003
$ws.Range("O7").Value = 0.4626735347223893
$ws.Range("P7").Value = 0.4626735347223893
$ws.Range("Q7").Value = 0.02711738440533333
$ws.Range("R7").Value = 0.244056459648
$ws.Range("S7").Value = 0.03633016256387022
$ws.Range("T7").Value = 0.03633016256387023
# Row 8
$ws.Range("A8").Value = 'Inflammatory-Mac'
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.06315533333333333
$ws.Range("H8").Value = 0.189466
$ws.Range("I8").Value = 0.07852224049441002
$ws.Range("J8").Value = 0.07852224049441003
$ws.Range("M8").Value = 0.01780266666666666
$ws.Range("N8").Value = 0.053408
$ws.Range("O8").Value = 0.01918323966442261
$ws.Range("P8").Value = 0.01918323966442261
$ws.Range("Q8").Value = 0.001124333347555555
$ws.Range("R8").Value = 0.010119000128
$ws.Range("S8").Value = 0.001506310958391697
$ws.Range("T8").Value = 0.001506310958391698
# Row 9
$ws.Range("A9").Value = 'Inflammatory-Mac'
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.06315533333333333
$ws.Range("H9").Value = 0.189466
$ws.Range("I9").Value = 0.07852224049441002
$ws.Range("J9").Value = 0.07852224049441003
$ws.Range("M9").Value = 0.3334263333333333
$ws.Range("N9").Value = 1.000279
$ws.Range("O9").Value = 0.3592830996908513
$ws.Range("P9").Value = 0.3592830996908513
$ws.Range("Q9").Value = 0.02105765122377777
$ws.Range("R9").Value = 0.189518861014
$ws.Range("S9").Value = 0.02821171395950211
$ws.Range("T9").Value = 0.02821171395950212
# Row 10
$ws.Range("A10").Value = 'MuSCs'
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4426526666666666
$ws.Range("H10").Value = 1.327958
$ws.Range("I10").Value = 0.5503585732663155
$ws.Range("J10").Value = 0.5503585732663157
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1474273333333333
$ws.Range("N10").Value = 0.442282
$ws.Range("O10").Value = 0.1588601259223368
$ws.Range("P10").Value = 0.1588601259223368
$ws.Range("Q10").Value = 0.06525910223955556
$ws.Range("R10").Value = 0.587331920156
$ws.Range("S10").Value = 0.08743003225152449
$ws.Range("T10").Value = 0.0874300322515245
# Row 11
$ws.Range("A11").Value = 'MuSCs'
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.4426526666666666
$ws.Range("H11").Value = 1.327958
$ws.Range("I11").Value = 0.5503585732663155
$ws.Range("J11").Value = 0.5503585732663157
$ws.Range("O11").Value = 0.4626735347223893
$ws.Range("P11").Value = 0.4626735347223893
$ws.Range("Q11").Value = 0.1900644314026667
$ws.Range("R11").Value = 1.710579882624
$ws.Range("S11").Value = 0.2546363464578973
$ws.Range("T11").Value = 0.2546363464578973
# Row 12
$ws.Range("A12").Value = 'MuSCs'
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4426526666666666
$ws.Range("H12").Value = 1.327958
$ws.Range("I12").Value = 0.5503585732663155
$ws.Range("J12").Value = 0.5503585732663157
$ws.Range("M12").Value = 0.01780266666666666
$ws.Range("N12").Value = 0.053408
$ws.Range("O12").Value = 0.01918323966442261
$ws.Range("P12").Value = 0.01918323966442261
$ws.Range("Q12").Value = 0.007880397873777777
$ws.Range("R12").Value = 0.070923580864
$ws.Range("S12").Value = 0.01055766041233742
$ws.Range("T12").Value = 0.01055766041233742
# Row 13
$ws.Range("A13").Value = 'MuSCs'
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4426526666666666
$ws.Range("H13").Value = 1.327958
$ws.Range("I13").Value = 0.5503585732663155
$ws.Range("J13").Value = 0.5503585732663157
$ws.Range("M13").Value = 0.3334263333333333
$ws.Range("N13").Value = 1.000279
$ws.Range("O13").Value = 0.3592830996908513
$ws.Range("P13").Value = 0.3592830996908513
$ws.Range("Q13").Value = 0.1475920555868889
$ws.Range("R13").Value = 1.328328500282
$ws.Range("S13").Value = 0.1977345341445563
$ws.Range("T13").Value = 0.1977345341445564
# Row 14
$ws.Range("G14").Value = 0.08209333333333334
$ws.Range("H14").Value = 0.24628
$ws.Range("I14").Value = 0.1020682200973436
$ws.Range("J14").Value = 0.1020682200973436
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1474273333333333
$ws.Range("N14").Value = 0.442282
$ws.Range("O14").Value = 0.1588601259223368
$ws.Range("P14").Value = 0.1588601259223368
$ws.Range("Q14").Value = 0.01210280121777778
$ws.Range("R14").Value = 0.10892521096
$ws.Range("S14").Value = 0.01621457029733279
$ws.Range("T14").Value = 0.01621457029733279
# Row 15
$ws.Range("G15").Value = 0.08209333333333334
$ws.Range("H15").Value = 0.24628
$ws.Range("I15").Value = 0.1020682200973436
$ws.Range("J15").Value = 0.1020682200973436
$ws.Range("O15").Value = 0.4626735347223893
$ws.Range("P15").Value = 0.4626735347223893
$ws.Range("Q15").Value = 0.03524890709333333
$ws.Range("R15").Value = 0.31724016384
$ws.Range("S15").Value = 0.04722426417526078
$ws.Range("T15").Value = 0.04722426417526078
# Row 16
$ws.Range("G16").Value = 0.08209333333333334
$ws.Range("H16").Value = 0.24628
$ws.Range("I16").Value = 0.1020682200973436
$ws.Range("J16").Value = 0.1020682200973436
$ws.Range("M16").Value = 0.01780266666666666
$ws.Range("N16").Value = 0.053408
$ws.Range("O16").Value = 0.01918323966442261
$ws.Range("P16").Value = 0.01918323966442261
$ws.Range("Q16").Value = 0.001461480248888889
$ws.Range("R16").Value = 0.01315332224
$ws.Range("S16").Value = 0.001957999128248379
$ws.Range("T16").Value = 0.001957999128248379
# Row 17
$ws.Range("G17").Value = 0.08209333333333334
$ws.Range("H17").Value = 0.24628
$ws.Range("I17").Value = 0.1020682200973436
$ws.Range("J17").Value = 0.1020682200973436
$ws.Range("M17").Value = 0.3334263333333333
$ws.Range("N17").Value = 1.000279
$ws.Range("O17").Value = 0.3592830996908513
$ws.Range("P17").Value = 0.3592830996908513
$ws.Range("Q17").Value = 0.02737207912444444
$ws.Range("R17").Value = 0.24634871212
$ws.Range("S17").Value = 0.03667138649650165
$ws.Range("T17").Value = 0.03667138649650165
